# Append 21 new weekly ECBASSETSW observations (rows 1272-1292) to Sheet1,
# continuing directly after the last existing data row (1271).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$symbol = ":ECBASSETSW"

# r, date-serial(A), open(C), high(D), low(E), close(F), volume(G)
$rows = @(
    @(1272, 44981, 7839385000000, 7839385000000, 7839385000000, 7839385000000, 0),
    @(1273, 44988, 7830796000000, 7830796000000, 7830796000000, 7830796000000, 0),
    @(1274, 44995, 7829368000000, 7829368000000, 7829368000000, 7829368000000, 0),
    @(1275, 45002, 7831859000000, 7831859000000, 7831859000000, 7831859000000, 0),
    @(1276, 45009, 7835921000000, 7835921000000, 7835921000000, 7835921000000, 0),
    @(1277, 45016, 7729612000000, 7729612000000, 7729612000000, 7729612000000, 0),
    @(1278, 45023, 7729251000000, 7729251000000, 7729251000000, 7729251000000, 0),
    @(1279, 45030, 7730924000000, 7730924000000, 7730924000000, 7730924000000, 0),
    @(1280, 45037, 7714222000000, 7714222000000, 7714222000000, 7714222000000, 0),
    @(1281, 45044, 7719816000000, 7719816000000, 7719816000000, 7719816000000, 0),
    @(1282, 45051, 7716913000000, 7716913000000, 7716913000000, 7716913000000, 0),
    @(1283, 45058, 7728510000000, 7728510000000, 7728510000000, 7728510000000, 0),
    @(1284, 45065, 7730118000000, 7730118000000, 7730118000000, 7730118000000, 0),
    @(1285, 45072, 7713658000000, 7713658000000, 7713658000000, 7713658000000, 0),
    @(1286, 45079, 7712715000000, 7712715000000, 7712715000000, 7712715000000, 0),
    @(1287, 45086, 7714391000000, 7714391000000, 7714391000000, 7714391000000, 0),
    @(1288, 45093, 7709739000000, 7709739000000, 7709739000000, 7709739000000, 0),
    @(1289, 45100, 7710607000000, 7710607000000, 7710607000000, 7710607000000, 0),
    @(1290, 45107, 7219695000000, 7219695000000, 7219695000000, 7219695000000, 0),
    @(1291, 45114, 7206902000000, 7206902000000, 7206902000000, 7206902000000, 0),
    @(1292, 45121, 7205494000000, 7205494000000, 7205494000000, 7205494000000, 0)
)

foreach ($row in $rows) {
    $r = $row[0]

    # Column A carries the date-formatted style (s="2") used throughout the
    # column; copy it down from the previous row before writing the value.
    $ws.Range("A1271").Copy($ws.Range("A$r"))
    $ws.Range("A$r").Value = $row[1]

    $ws.Range("B$r").Value = $symbol

    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}
